# Insert a new weekly record for "Jengibre" at Terminal La Palmera de La Serena.
# This shifts the existing rows 116..164 down to 117..165 and inserts a brand
# new row 116 with the latest week's data (the rest of the table is otherwise
# unchanged, just pushed down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 116; Excel shifts rows 116-164 down to 117-165
# and copies formatting (incl. the date number format on column D) from the
# row above automatically.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new record.
$ws.Cells.Item(116, 1).Value = 8
$ws.Cells.Item(116, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 45141
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = 100114007
$ws.Cells.Item(116, 7).Value = "Jengibre"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 300
$ws.Cells.Item(116, 11).Value = 18000
$ws.Cells.Item(116, 12).Value = 19000
$ws.Cells.Item(116, 13).Value = 18500
$ws.Cells.Item(116, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(116, 15).Value = "Perú"
$ws.Cells.Item(116, 16).Value = 1423
$ws.Cells.Item(116, 17).Value = 13
$ws.Cells.Item(116, 18).Value = "Hortaliza"
